# "changed signin to be the first step"
# Insert a new product row (ASUS TUF Gaming NVIDIA GeForce RTX 3070 OC Edition)
# right after the "ASUS ROG STRIX" row, pushing the following rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before current row 5 (EVGA row), shifting rows 5-9 down to 6-10.
$ws.Rows.Item(5).Insert()

# Fill in the new row's data: ASIN (B) first, then Description (A), then MSRP (C)
# so the shared-string table gets the ASIN before the description (matches
# how the workbook author entered the data).
$ws.Cells.Item(5, 2).Value = "B08L8KC1J7"
$ws.Cells.Item(5, 1).Value = "ASUS TUF Gaming NVIDIA GeForce RTX 3070 OC Edition "
$ws.Cells.Item(5, 3).Value = 800

# Update the selection / active cell like the resaved workbook shows.
$ws.Range("C12").Select()
